$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4

$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 0.06

$ws.Range("B7").Value = 11

$ws.Range("B8").Value = 16
$ws.Range("C8").Value = 0.04

$ws.Range("C11").Value = 0.01

$ws.Range("B13").Value = 22
$ws.Range("C13").Value = 0.06

$ws.Range("C14").Value = 0.02

$ws.Range("B15").Value = 32
$ws.Range("C15").Value = 0.09

$ws.Range("B16").Value = 9

$ws.Range("B17").Value = 14

$ws.Range("B18").Value = 9
$ws.Range("C18").Value = 0.03

$ws.Range("B19").Value = 11

$ws.Range("B22").Value = 30
$ws.Range("C22").Value = 0.08

$ws.Range("B24").Value = 8

$ws.Range("C26").Value = 0.04

$ws.Range("B27").Value = 3
$ws.Range("C27").Value = 0.01

$ws.Range("B29").Value = 71

$ws.Range("B31").Value = 12
$ws.Range("C31").Value = 0.03
